$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 469.8
$ws.Range("I2").Value = 499
$ws.Range("J2").Value = 462.5
$ws.Range("K2").Value = 499
$ws.Range("L2").Value = 462.5
$ws.Range("M2").Value = -386
$ws.Range("N2").Value = -688.5
$ws.Range("H6").Value = 680.4375
$ws.Range("I6").Value = 353.36365
$ws.Range("J6").Value = 1400
$ws.Range("K6").Value = 1060.09095
$ws.Range("L6").Value = 4200
$ws.Range("M6").Value = -948.09095
$ws.Range("N6").Value = -4424
$ws.Range("H21").Value = 11000
$ws.Range("J21").Value = 11000
$ws.Range("L21").Value = 11000
$ws.Range("N21").Value = -11936
$ws.Range("H23").Value = 11000
$ws.Range("J23").Value = 11000
$ws.Range("L23").Value = 11000
$ws.Range("N23").Value = -11468
$ws.Range("H40").Value = 2207.625
$ws.Range("I40").Value = 2281.4285
$ws.Range("J40").Value = 2150.2222
$ws.Range("K40").Value = 2281.4285
$ws.Range("L40").Value = 2150.2222
$ws.Range("M40").Value = -2106.4285
$ws.Range("N40").Value = -2500.2222
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 36903.387
$ws.Range("I132").Value = 37256.895
$ws.Range("J132").Value = 33604
$ws.Range("K132").Value = 111770.685
$ws.Range("L132").Value = 100812
$ws.Range("M132").Value = -109240.685
$ws.Range("N132").Value = -105872
$ws.Range("H141").Value = 3851.037
$ws.Range("I141").Value = 2693.238
$ws.Range("K141").Value = 8079.714
$ws.Range("M141").Value = -2899.714

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36339.766
$ws.Range("I32").Value = 5751.2666
$ws.Range("J32").Value = 265753.5
$ws.Range("K32").Value = 5751.2666
$ws.Range("L32").Value = 265753.5
$ws.Range("M32").Value = -5464.2666
$ws.Range("N32").Value = -266327.5
$ws.Range("H61").Value = 3383.5715
$ws.Range("I61").Value = 2901.4443
$ws.Range("J61").Value = 5010.75
$ws.Range("K61").Value = 2901.4443
$ws.Range("L61").Value = 5010.75
$ws.Range("M61").Value = -2689.4443
$ws.Range("N61").Value = -5434.75
$ws.Range("H122").Value = 4775
$ws.Range("I122").Value = 2550
$ws.Range("K122").Value = 7650
$ws.Range("M122").Value = -5200
$ws.Range("H136").Value = 3383.5715
$ws.Range("I136").Value = 2901.4443
$ws.Range("J136").Value = 5010.75
$ws.Range("K136").Value = 8704.332900000001
$ws.Range("L136").Value = 15032.25
$ws.Range("M136").Value = -6154.332900000001
$ws.Range("N136").Value = -20132.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3541.6453
$ws.Range("I105").Value = 3340.5264
$ws.Range("J105").Value = 3860.0833
$ws.Range("K105").Value = 3340.5264
$ws.Range("L105").Value = 3860.0833
$ws.Range("M105").Value = -1593.5264
$ws.Range("N105").Value = -7354.0833
$ws.Range("H122").Value = 30153.166
$ws.Range("J122").Value = 30153.166
$ws.Range("L122").Value = 30153.166
$ws.Range("N122").Value = -39953.166

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1194.6
$ws.Range("I16").Value = 980
$ws.Range("K16").Value = 980
$ws.Range("M16").Value = -693
$ws.Range("H31").Value = 5402.7295
$ws.Range("I31").Value = 2000.091
$ws.Range("J31").Value = 10393.267
$ws.Range("K31").Value = 2000.091
$ws.Range("L31").Value = 10393.267
$ws.Range("M31").Value = -1705.091
$ws.Range("N31").Value = -10983.267
$ws.Range("H34").Value = 5402.7295
$ws.Range("I34").Value = 2000.091
$ws.Range("J34").Value = 10393.267
$ws.Range("K34").Value = 2000.091
$ws.Range("L34").Value = 10393.267
$ws.Range("M34").Value = -1798.091
$ws.Range("N34").Value = -10797.267
$ws.Range("H113").Value = 1194.6
$ws.Range("I113").Value = 980
$ws.Range("K113").Value = 980
$ws.Range("M113").Value = 1190

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 62500164
$ws.Range("J2").Value = 200000460
$ws.Range("L2").Value = 1200002760
$ws.Range("N2").Value = -1200002986
$ws.Range("H38").Value = 124.4
$ws.Range("I38").Value = 142.33333
$ws.Range("J38").Value = 116.71429
$ws.Range("K38").Value = 426.99999
$ws.Range("L38").Value = 350.14287
$ws.Range("M38").Value = -79.99998999999997
$ws.Range("N38").Value = -1044.14287
$ws.Range("H68").Value = 20506
$ws.Range("J68").Value = 452.5
$ws.Range("L68").Value = 1357.5
$ws.Range("N68").Value = -2979.5
$ws.Range("H71").Value = 20506
$ws.Range("J71").Value = 452.5
$ws.Range("L71").Value = 4072.5
$ws.Range("N71").Value = -12184.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7489.231
$ws.Range("I70").Value = 9296.666999999999
$ws.Range("J70").Value = 5940
$ws.Range("K70").Value = 9296.666999999999
$ws.Range("L70").Value = 5940
$ws.Range("M70").Value = -9026.666999999999
$ws.Range("N70").Value = -6480
$ws.Range("H73").Value = 7489.231
$ws.Range("I73").Value = 9296.666999999999
$ws.Range("J73").Value = 5940
$ws.Range("K73").Value = 9296.666999999999
$ws.Range("L73").Value = 5940
$ws.Range("M73").Value = -8360.666999999999
$ws.Range("N73").Value = -7812
$ws.Range("H123").Value = 15230.615
$ws.Range("I123").Value = 6350
$ws.Range("J123").Value = 16845.273
$ws.Range("K123").Value = 6350
$ws.Range("L123").Value = 16845.273
$ws.Range("M123").Value = -3900
$ws.Range("N123").Value = -21745.273

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2933.2173
$ws.Range("I40").Value = 1692.8
$ws.Range("J40").Value = 3277.7778
$ws.Range("K40").Value = 1692.8
$ws.Range("L40").Value = 3277.7778
$ws.Range("M40").Value = -1556.8
$ws.Range("N40").Value = -3549.7778
$ws.Range("H122").Value = 3391.318
$ws.Range("I122").Value = 2229.75
$ws.Range("J122").Value = 3649.4443
$ws.Range("K122").Value = 6689.25
$ws.Range("L122").Value = 10948.3329
$ws.Range("M122").Value = -4239.25
$ws.Range("N122").Value = -15848.3329

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 27559.143
$ws.Range("I64").Value = 28000
$ws.Range("J64").Value = 27485.666
$ws.Range("K64").Value = 28000
$ws.Range("L64").Value = 27485.666
$ws.Range("M64").Value = -27752
$ws.Range("N64").Value = -27981.666
$ws.Range("H67").Value = 27559.143
$ws.Range("I67").Value = 28000
$ws.Range("J67").Value = 27485.666
$ws.Range("K67").Value = 28000
$ws.Range("L67").Value = 27485.666
$ws.Range("M67").Value = -27142
$ws.Range("N67").Value = -29201.666
